$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.855.64"
$ws.Range("E2").Value = "  +1.24%  "
$ws.Range("D3").Value = "3.966.69"
$ws.Range("E3").Value = "  +4.04%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "486.50"
$ws.Range("E5").Value = "  +8.71%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "149.52"
$ws.Range("E6").Value = "  +3.69%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.629"
$ws.Range("E7").Value = "  +1.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.998"
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.733"
$ws.Range("E9").Value = "  -0.67%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.174"
$ws.Range("E10").Value = "  +13.68%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000378"
$ws.Range("E11").Value = "  +20.77%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "43.56"
$ws.Range("E12").Value = "  -0.25%  "
$ws.Range("D13").Value = "4.600.62"
$ws.Range("E13").Value = "  +4.43%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.43"
$ws.Range("E14").Value = "  +0.44%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.92"
$ws.Range("E15").Value = "  +0.72%  "
$ws.Range("D16").Value = "3.974.24"
$ws.Range("E16").Value = "  +2.62%  "
$ws.Range("E17").Value = "  +0.32%  "
$ws.Range("E18").Value = "  +0.18%  "
$ws.Range("E19").Value = "  +0.69%  "
$ws.Range("D20").Value = "67.958.15"
$ws.Range("E20").Value = "  +1.31%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "434.47"
$ws.Range("E21").Value = "  +3.08%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.41"
$ws.Range("E22").Value = "  +3.80%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.45"
$ws.Range("E23").Value = "  -1.23%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "88.21"
$ws.Range("E24").Value = "  +2.59%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.65"
$ws.Range("E25").Value = "  +5.91%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "38.84"
$ws.Range("E26").Value = "  +4.19%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.06"
$ws.Range("E27").Value = "  +3.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.74"
$ws.Range("E28").Value = "  +3.70%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "721.52"
$ws.Range("E29").Value = "  -1.34%  "
$ws.Range("E30").Value = "  -2.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "13.35"
$ws.Range("E31").Value = "  -3.25%  "
$ws.Range("E32").Value = "  +2.65%  "
$ws.Range("B33").Value = "PEPE"
$ws.Range("C33").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D33").Value = "0.0₃0897"
$ws.Range("E33").Value = "  +32.48%  "
$ws.Range("B34").Value = "InjectiveProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "41.99"
$ws.Range("E34").Value = "  -4.35%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "60.70"
$ws.Range("E35").Value = "  +7.30%  "
$ws.Range("E36").Value = "  -4.18%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.998"
$ws.Range("E37").Value = "  -0.18%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.36"
$ws.Range("E38").Value = "  -3.18%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0475"
$ws.Range("E39").Value = "  -1.33%  "
$ws.Range("E40").Value = "  +5.39%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.143"
$ws.Range("E41").Value = "  +1.74%  "
$ws.Range("B42").Value = "WEMIXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.87"
$ws.Range("E42").Value = "  +8.89%  "
$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.24"
$ws.Range("E43").Value = "  +4.94%  "
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("E44").Value = "  -0.08%  "
$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.335"
$ws.Range("E45").Value = "  -1.66%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.44"
$ws.Range("E46").Value = "  +2.06%  "
$ws.Range("E47").Value = "  -2.86%  "
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.23"
$ws.Range("E48").Value = "  -0.91%  "
$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "148.96"
$ws.Range("E49").Value = "  +2.44%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.91"
$ws.Range("E50").Value = "  +1.40%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "25.06"
$ws.Range("E51").Value = "  +0.09%  "
